$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("List1")

# Numeric "Hodiny"/weight values in column F that were updated
$ws.Range("F2").Value = 3
$ws.Range("F16").Value = 1.2
$ws.Range("F17").Value = 2
$ws.Range("F18").Value = 2.3
$ws.Range("F20").Value = 1.2
$ws.Range("F21").Value = 1.2
$ws.Range("F22").Value = 2
$ws.Range("F23").Value = 2.3
$ws.Range("F24").Value = 2.3
$ws.Range("F25").Value = 1.2

# Descriptive comment text in column G that was rewritten
$ws.Range("G22").Value = "vzorců tam je hodne ale dá se to"
$ws.Range("G23").Value = "extrém, ale dá se"
